$wb = $excel.ActiveWorkbook

# ---- Sheet YDS: append Week 13 play-by-play yardage logs ----
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = "14 2 -6 0 9 7 3 9 3 16 3 1 3 1 25 2 8 2 6 5 2 7 0 5 6 9 7 10 80 -2 0 6 -2 -3 1 -2 7 0 -3 5 2 6 55 -4 -1 3 5 -2 5 16 0 1 1 2 3 1 9 0 3 4 4 -3 -1 5 4 4 0 -1 -3 11 10 5 2 0 19 19 1 -1 2 3 1 23 3 6 3 -2 2 5 15 10 8 -5 38 3 4 -1 5 3 3 0 0 9 4 11 2 3 1 0 1 5 0 7 5 5 37 5 7 11 4 0 9 4 3 7 9 9 -1 5 5 2 2 7 4 3 15 3 4 0 3 8 5 2 2 1 3 6 0 1 1 5 -6 4 1 9 4 4 4 2 1 6 10 0 3 3 4 4 2 3 -2 5 3 13 9 1 4 7 17 6 3 4 20 5 16 8 5 7 3 6 6 9 10 -1 1 2 0 9 3 2 1 -1 4 5 -3 4 12 3 3 2 3 2 -1 3 1 1 3 4 1 14 2 1 3 8 3 -1 6 2 11 4 1 -1 8 2 0 3 5 1 -3 2 6 0 1 7 -8 -4 5 5 3 -2 1 8 4 3 2 -2 9 6 1 2 2 3 4 6 7 3 5 3 8 0 2 10 2 6 10 19 0 5 4 4 -4 1 7 5 -2 0 0 0 1 3 1 2 0 4 5 -2 16 8 4 4 8 1 0 0 1 3 2 12 3 18 -9 4 10 9 2 5 3 9 9 4 11 1 1 6 1 0 6 -2 2 -2 6 4 8 1 14 3 0 9 0 1 1 8 2 5 1 8 4 9 17 -3 2 8 10 5 7 -1 0 1 1 9 0 5 10 1 6 3 3 1 0 3 1 4 12 -1 1 7 8 6 3 3 6 3 7 14 17 11 3 1 4 7 15 16 5 2 23 2 0 -3 2 34 25 -2 12 0 4 6 0 7 7 0 0 3 7 8 2 -1 3 18 -4 4 3 -1 2 1 5 7 4 4 1 0 11 9 1 3 3 38 6 11 3 3 3 2 10 2 7 3 2 1 7 -1 0 0 6 0 1 1 -4 -1 3 3 0 1 2 5 2 7 4 5 2 2 10 2 3 2 7 4 21 2 0 1 3 9 3 0 8 -1 8 1 1 3 4 -3 0 9 2 2 0 1 2 2 16 8 4 0 7 9 5 1 1 3 -1 5 6 9 2 7 4 2 8 4 5 5 -1 5 15 7 -1 1 6 7 2 4 3 3 6 1 7 7 14 0 2 1 -3 5 11 9 9 7 7 14 4 -2 1 7 7 12 1 13 -2 15 1 5 6 7 4 5 3 20 14 14 12 6 2 3 6 5 3 -1 2 3 -3 2 0 4 1 5 4 7 1 1 27 1 -1 3 0 4 2 3 11 5 27 1 5 39 -1 12 -2 0 1 9 -5 -2 11 14 5 2 -1 2 3 4 7 5 1 9 4 0 2 0 8 6 5 -1 9 -1 1 6 2 13 2 8 17 4 -2 4 4 4 -6 3 -3 4 8 4 2 0 -2 2 1 9 2 0 -4 8 4 4 1 2 2 2 4 0 12 3 2 9 0 1 3 11 -1 25 3 5 -1 1 0 4 4 3 11 21 -2 12 0 5 2 4 4 -1 3 2 3 7 4 3 -3 6 15 20 10 1 3 4 9 3 4 10 4 4 2 5 49 4 3 2 4 12 0 1 3 -6 5 1 0 6 7 9 4 3 1 1 -1 4 5 0 6 2 7 5 2 9 2 4 0 2 2 5 4 0 2 6 3 1"
$wsYDS.Range("B3").Value = "5 11 76 5 14 14 10 18 9 6 2 3 9 12 16 41 5 5 6 2 9 5 4 18 8 4 11 10 4 18 18 8 9 5 10 4 10 15 5 1 24 13 2 1 20 14 11 12 9 1 14 7 14 17 20 14 26 8 14 19 6 5 11 19 19 18 9 20 16 9 17 7 13 0 17 38 13 5 4 10 12 13 2 5 20 5 10 8 7 5 11 13 4 14 14 2 20 13 5 10 12 15 5 11 15 11 20 0 12 19 15 12 5 9 9 4 9 35 10 18 6 2 7 4 8 5 19 25 8 2 11 4 7 -2 8 18 7 11 11 7 23 15 14 1 14 19 20 7 8 15 16 9 13 23 6 35 18 3 11 19 8 9 13 6 2 10 14 3 -6 6 12 9 25 12 4 9 9 33 2 20 11 14 17 29 7 3 6 3 43 4 5 18 3 5 6 1 47 10 9 36 41 8 11 7 9 4 8 7 12 5 8 4 14 4 14 22 13 8 7 12 11 13 8 7 1 11 -1 -5 27 4 26 22 4 33 26 3 4 12 9 3 -6 18 4 26 -1 1 29 4 6 3 12 20 7 24 7 4 17 11 2 13 13 1 17 20 7 4 12 49 6 22 15 9 18 2 3 10 12 23 5 11 9 5 16 6 6 5 17 18 4 11 8 12 12 16 10 -1 13 5 16 6 9 6 12 14 7 19 14 12 13 5 22 19 9 -3 4 23 25 3 13 3 9 15 5 11 7 8 -1 23 3 26 49 15 24 21 16 11 18 27 0 10 9 7 1 23 3 3 1 -5 45 4 14 2 4 9 3 6 41 27 5 13 12 5 5 19 21 8 12 14 3 23 7 5 29 8 7 23 7 14 18 35 8 79 16 2 15 12 11 10 1 8 7 6 11 17 5 8 5 6 2 40 11 9 4 6 6 6 0 15 8 8 7 6 4 16 8 16 10 11 5 7 19 1 0 16 8 0 16 12 9 3 39 12 12 12 9 8 21 9 18 21 7 15 14 26 7 10 8 8 14 10 12 12 7 8 16 14 26 6 1 18 6 12 4 26 5 13 26 13 9 9 0 28 10 10 16 6 1 29 27 14 3 9 5 5 21 -4 19 11 6 23 16 6 15 50 3 83 15 7 11 19 23 4 19 18 16 7 11 14 22 11 10 16 6 5 19 29 1 11 32 7 10 15 6 5 12 5 3 7 5 -1 18 6 7 8 10 20 21 6 3 19 7 11 7 9 7 8 13 34 7 6 9 9 23 15 7 13 4 11 21 37 3 2 11 2 9 24 3 13 13 30 12 8 11 5 26 24 8 28 12 4 7 9 8 9 48 3 21 13 29 11 18 14 10 7 16"
$wsYDS.Range("C2").Value = "2 1 4 3 11 0 9 5 1 20 0 11 7 3 -4 2 3 8 -1 13 2 2 2 22 10 1 15 1 9 2 5 7 0 5 1 -2 8 4 -2 3 13 2 2 4 5 0 4 3 1 12 2 5 -1 3 10 6 5 0 4 0 6 19 3 2 1 23 -3 7 3 16 1 0 2 0 6 0 1 8 5 8 11 3 -1 8 9 0 14 4 8 7 2 2 5 -3 4 12 0 1 0 -1 6 3 4 1 -2 2 17 5 -1 4 4 3 0 12 2 5 4 21 2 12 1 1 0 0 2 9 0 0 2 2 1 10 18 4 7 11 5 3 1 11 0 13 9 4 7 -1 6 5 2 4 5 6 2 2 3 7 12 2 4 3 3 9 10 6 -3 2 2 2 6 -3 5 -1 13 4 3 5 -1 21 9 4 2 1 0 6 15 4 0 0 1 4 4 1 6 1 11 6 2 2 9 2 11 6 3 2 5 3 7 2 3 0 4 4 2 3 8 4 -2 9 0 0 0 2 0 4 5 -6 5 0 0 7 2 5 3 1 3 8 6 1 2 18 3 1 5 9 11 3 10 6 3 7 2 -3 3 2 4 2 -1 4 5 8 -1 1 5 3 1 3 2 6 3 -1 1 61 6 1 4 0 6 4 5 7 3 0 -2 4 0 2 5 4 0 8 0 2 6 1 8 15 6 -1 2 2 6 -2 3 11 3 5 3 3 7 -2 2 3 9 17 9 0 17 -1 5 3 5 3 1 0 2 2 3 3 2 1 2 7 2 0 5 3 -1 2 6 1 4 1 2 4 5 11 6 40 7 6 10 6 2 4 3 1 11 3 2 2 -1 8 -3 0 0 0 1 5 2 0 -1 -2 3 3 7 2 0 1 4 5 2 2 2 1 9 1 9 5 11 3 16 9 1 6 6 8 8 10 2 4 9 20 7 7 1 1 3 16 15 -1 -4 3 -4 1 3 6 1 10 1 2 2 0 5 10 5 2 6 6 12 -2 9 -2 -3 1 7 9 7 6 6 13 3 3 0 8 9 27 1 2 5 -4 1 3 2 4 2 10 4 2 6 1 6 3 7 2 6 2 12 4 2 12 4 2 2 4 2 -2 -1 5 4 3 5 2 12 -1 6 8 6 16 1 6 4 3 -2 4 14 0 1 5 2 4 5 26 1 3 9 1 5 0 12 1 0 -2 11 4 1 2 10 1 2 8 4 2 6 2 0 33 -2 4 16 1 1 5 10 -2 0 1 9 8 3 3 0 5 0 8 5 17 -1 0 5 0 5 3 -1 -1 4 22 9 -2 1 3 6 9 6 3 2 3 4 16 3 9 3 16 2 0 15 -3 -4 4 -4 4 6 3 0 14 8 4 0 0 3 8 1 4 1 4 2 1 1 2 3 3 1 -4 1 5 8 2 35 21 5 1 7 6 2 4 1 3 5 2 3 -3 2 6 -2 10 6 1 3 11 4 3 6 7 -1 1 7 10 -3 2 7 2 2 4 4 4 4 4 6 1 4 1 3 1 1 0 4 7 4 30 0 6 -4 7 1 3 1 2 -5 0 11 12 0 7 4 2 1 4 2 3 3 3 3 2 0 4 1 10 4 2 0 -1"
$wsYDS.Range("C3").Value = "3 7 10 1 2 10 13 6 9 12 15 5 13 12 5 8 7 18 7 11 0 4 33 9 9 7 3 8 8 5 5 11 1 18 1 13 5 0 27 2 4 8 5 30 2 18 9 13 23 0 11 2 9 8 5 5 10 29 4 22 -4 6 20 28 7 8 4 5 15 13 23 9 3 1 2 2 17 17 17 42 47 15 3 28 6 6 8 12 19 2 22 70 4 1 8 4 7 6 15 31 32 4 17 3 10 8 -5 31 12 6 17 7 4 4 7 6 8 5 5 4 40 6 3 1 20 21 12 16 11 8 8 3 13 17 8 15 -3 7 5 7 3 9 12 46 11 15 35 6 2 4 3 7 2 9 6 14 12 5 8 17 6 9 4 5 15 36 11 5 4 7 1 2 11 52 49 3 11 1 8 1 24 9 7 3 34 6 3 15 1 8 34 7 5 -1 15 11 3 19 12 15 3 11 7 11 4 30 6 5 8 6 10 9 8 5 33 14 12 8 6 5 3 31 5 5 6 12 35 22 5 8 9 8 8 12 4 7 18 20 18 9 12 4 6 11 23 3 10 12 3 28 11 3 2 11 13 5 9 9 5 13 12 8 -1 9 5 5 5 12 4 10 3 30 1 11 4 13 8 22 1 12 2 45 12 9 5 12 6 38 5 -1 3 7 5 7 7 -1 8 3 5 16 5 5 -3 11 9 45 2 5 5 5 8 7 9 9 4 9 4 5 5 6 8 2 29 26 5 4 4 7 10 9 5 20 6 12 9 5 22 6 -8 12 12 1 3 6 2 6 7 10 -1 7 11 6 43 3 23 5 5 3 10 5 6 5 3 8 7 9 11 13 3 20 2 20 25 4 11 4 9 1 5 26 91 6 10 9 14 42 11 -4 2 12 2 47 14 1 0 10 3 4 15 19 6 1 9 6 7 12 25 17 -2 5 3 2 28 28 12 1 19 9 2 13 16 6 2 5 7 8 13 10 9 23 33 3 6 2 9 14 8 5 8 13 7 9 2 30 9 11 -1 57 11 -2 2 -1 -1 5 5 -2 2 14 13 9 5 6 28 14 6 11 11 16 6 8 10 4 2 -6 19 22 6 2 9 10 20 5 10 0 9 9 8 0 21 17 33 50 17 6 -2 25 45 4 4 1 5 5 9 6 0 5 2 6 8 16 29 13 10 8 33 5 -2 5 2 6 -3 12 16 6 11 5 9 5 2 25 9 6 6 18 16 6 8 15 7 3 7 11 13 22 8 9 2 11 29 7 4 30 15 20 7 6 14 15 24 16 5 17 7 5 23 2 -1 4 1 0 8 10 7 12 -1 33 5 9 2 4 11 7 6 7 16 17 5 7 2 12 11 5 2 27 3 0"

# ---- Sheet OFF: update Home (row2) and Road (row3) season totals ----
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 395
$wsOFF.Range("F2").Value = 106
$wsOFF.Range("G2").Value = 114
$wsOFF.Range("J2").Value = 45
$wsOFF.Range("N2").Value = 24
$wsOFF.Range("O2").Value = 32
$wsOFF.Range("P2").Value = 15
$wsOFF.Range("B3").Value = 23
$wsOFF.Range("C3").Value = 265
$wsOFF.Range("E3").Value = 78
$wsOFF.Range("F3").Value = 163
$wsOFF.Range("G3").Value = 51
$wsOFF.Range("H3").Value = 47
$wsOFF.Range("I3").Value = 98
$wsOFF.Range("J3").Value = 74
$wsOFF.Range("L3").Value = 430
$wsOFF.Range("M3").Value = 290
$wsOFF.Range("Q3").Value = 876

# ---- Sheet DEF: update Home (row2) and Road (row3) season totals ----
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 341
$wsDEF.Range("D2").Value = 18
$wsDEF.Range("E2").Value = 24
$wsDEF.Range("F2").Value = 103
$wsDEF.Range("G2").Value = 98
$wsDEF.Range("H2").Value = 10
$wsDEF.Range("J2").Value = 47
$wsDEF.Range("N2").Value = 35
$wsDEF.Range("O2").Value = 39
$wsDEF.Range("P2").Value = 21
$wsDEF.Range("B3").Value = 13
$wsDEF.Range("C3").Value = 299
$wsDEF.Range("E3").Value = 52
$wsDEF.Range("F3").Value = 166
$wsDEF.Range("G3").Value = 60
$wsDEF.Range("H3").Value = 39
$wsDEF.Range("I3").Value = 110
$wsDEF.Range("J3").Value = 100
$wsDEF.Range("L3").Value = 466
$wsDEF.Range("M3").Value = 302
$wsDEF.Range("Q3").Value = 850

# ---- Sheet ST: update special teams totals + append Week 13 kick/punt logs ----
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 140
$wsST.Range("D2").Value = 100
$wsST.Range("F2").Value = 574
$wsST.Range("G2").Value = 559
$wsST.Range("N2").Value = 49
$wsST.Range("O2").Value = 34
$wsST.Range("B4").Value = "70 70 65 66 61 69 43 65 60 59 58 61 66 66 62 67 66 52 69 63 67 56 42 69 65 64 67 62 66 63 69 69 67 64 56 65 63 62 64 52 70 54 68 73 63 61 60 66 63 58 66 57 61 60 67 70 64 51 55 65 64 64 63 66 61 56 69"
$wsST.Range("B5").Value = "42 42 27 17 32 18 0 18 16 15 15 17 30 22 21 29 31 15 31 26 21 17 5 26 75 21 25 30 20 25 33 21 27 19 0 20 21 33 18 14 31 10 34 24 30 27 21 29 24 22 26 17 17 16 23 28 21 18 30 21 99 22 23 25 33 26 19"
$wsST.Range("B6").Value = "16 16 29 26 20 21 21 18 20 18 20 16 30 18 34 12 20 21 14 13 23 15 19 15 16 11 18 16 19 16 19 11 68 16 20 26 7 21 10 22 27 27 33 29 22 23 37 0 21 12 12 23"
$wsST.Range("D3").Value = "46 30 51 51 47 47 46 59 51 47 39 40 43 32 49 44 51 48 43 45 39 51 56 56 60 58 44 57 53 50 47 46 54 49 60 55 47 50 42 54 36 30 47 41 44 45 54 48 44 58 48 47 46 45 44 31 50 26 44 49 40 49 45 44 43 33 48 52 49 35 42 59 56 41 57 45 54 37 45 33 58 39 42 65 48 38 34 44 57 28 41 46 40 45 65 51 67 60 37 40"
$wsST.Range("D4").Value = "0 8 0 0 0 0 0 0 0 8 0 0 6 5 15 0 9 5 0 0 0 0 15 9 0 0 0 0 0 0 0 0 0 7 0 0 4 -2 0 22 0 0 22 0 0 0 14 0 0 0 10 0 16 0 0 0 0 0 0 0 0 7 7 0 0 0 11 2 0 0 0 0 0 0 6 0 0 0 11 0 12 0 5 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 11"
$wsST.Range("D5").Value = "0 0"

# ---- Sheet TURNS: update Road (row3) turnover totals ----
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 12
$wsTURNS.Range("C3").Value = 11
$wsTURNS.Range("E3").Value = 18

# ---- Sheet PEN: update penalty counts ----
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 31
$wsPEN.Range("B3").Value = 33
$wsPEN.Range("D4").Value = 23

# ---- Active sheet/selection: ST sheet active tab, D6 selected ----
$wsST.Activate()
$wsST.Range("D6").Select()

